$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the special date-only number format currently applied to A51
# (the last row's date cell) before we touch anything.
$lastRowFormat = $ws.Range("A51").NumberFormat

# Row 51 is no longer the last row, so its date cell switches to the same
# date/time format used by all the other non-final rows (e.g. A2).
$ws.Range("A51").NumberFormat = $ws.Range("A2").NumberFormat

# Append the new row 52 with the next day's data (8 AM UTC daily update).
$ws.Range("A52").Value = 45637
$ws.Range("B52").Value = 129
$ws.Range("C52").Value = 115
$ws.Range("D52").Value = 119

# The new last row's date cell takes on the special "last row" format
# that A51 used to have.
$ws.Range("A52").NumberFormat = $lastRowFormat
